# Insert a new row before row 547 (shifts existing rows 547..612 down to 548..613)
# and populate the new row 547 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 547
$ws.Rows.Item(547).Insert()

# Populate the newly inserted row 547 with the new data record.
$ws.Cells.Item(547, 1).Value = 5
$ws.Cells.Item(547, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(547, 3).Value = 'Maule'
$ws.Cells.Item(547, 4).Value = 45212
$ws.Cells.Item(547, 5).Value = 7
$ws.Cells.Item(547, 6).Value = 100114013
$ws.Cells.Item(547, 7).Value = 'Zanahoria'
$ws.Cells.Item(547, 8).Value = 'Sin especificar'
$ws.Cells.Item(547, 9).Value = 'Primera'
$ws.Cells.Item(547, 10).Value = 500
$ws.Cells.Item(547, 11).Value = 5000
$ws.Cells.Item(547, 12).Value = 5000
$ws.Cells.Item(547, 13).Value = 5000
$ws.Cells.Item(547, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(547, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(547, 16).Value = 250
$ws.Cells.Item(547, 17).Value = 20
$ws.Cells.Item(547, 18).Value = 'Hortaliza'

# Apply the date number format (matching column D's style elsewhere) to the new D547 cell.
$ws.Cells.Item(547, 4).NumberFormat = $ws.Cells.Item(548, 4).NumberFormat
